$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2,3) {
    $ws.Range("K$row").Value = -0.59
    $ws.Range("U$row").Value = 1.24
    $ws.Range("V$row").Value = 0.1024793388429752
    $ws.Range("W$row").Value = -0.05728155339805825
    $ws.Range("X$row").Value = 0.05024160799847834
    $ws.Range("Y$row").Value = -0.1075231613965366
    $ws.Range("AA$row").Value = -0.07175453040648792
    $ws.Range("AB$row").Value = 0.05021404256626861
    $ws.Range("AC$row").Value = -0.1219685729727565
    $ws.Range("AE$row").Value = 0.01165071166430416
    $ws.Range("AF$row").Value = 0.01165071166430416
    $ws.Range("AG$row").Value = -1.228349288335696
    $ws.Range("AH$row").Value = 0.0009619425082234062
    $ws.Range("AI$row").Value = 0.001129859029372004
    $ws.Range("AJ$row").Value = -0.1129864563269852
    $ws.Range("AK$row").Value = -0.1354052671754974
    $ws.Range("AM$row").Value = -0.002
    $ws.Range("AP$row").Value = 2.225270449883507
    $ws.Range("AQ$row").Value = 295.5
}
